# modifying an event in the calendar (drag, stretch, etc) updates the stats card.
#
# Mark the following backlog items as "done":
#   - Refactor controller (architecture)                                  (row 9,  C9)
#   - Either the work day is 8 hours, or the 9 = 100% (fix the status bar) (row 12, C12)
#   - When no event is selected ... stats card behaves correctly          (row 14, C14)
#
# The sheet has an AutoFilter on the Status column (col C) showing only
# "pending" rows. Re-apply it over the full data range (A1:C14, which now
# includes rows 10-14) so the newly-"done" rows get hidden and the filter
# range/FilterDatabase name/selection are kept in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the three statuses to "done".
$ws.Range("C9").Value = "done"
$ws.Range("C12").Value = "done"
$ws.Range("C14").Value = "done"

# 2. Re-apply the AutoFilter over the expanded range A1:C14 so rows that no
#    longer match "pending" become hidden, and the autoFilter ref grows.
$ws.AutoFilterMode = $false
$ws.Range("A1:C14").AutoFilter(3, "pending", 7)

# 3. Keep the hidden _xlnm._FilterDatabase defined name in sync with the
#    new filter range.
$names = $wb.Names
$filterDatabaseName = $names.Item(1)
$filterDatabaseName.RefersTo = "=Sheet1!`$A`$1:`$C`$14"

# 4. Move the active selection to C15 (just below the table).
$ws.Range("C15").Select()
